$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 09:19:21"
$wsZhCn.Range("E3").Value = "2016-03-11 09:19:21"
$wsZhCn.Range("H2").Value = "2016-03-11 09:19:38"
$wsZhCn.Range("H3").Value = "2016-03-11 09:19:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 09:19:25"
$wsDeDe.Range("E3").Value = "2016-03-11 09:19:25"
$wsDeDe.Range("H2").Value = "2016-03-11 09:19:44"
$wsDeDe.Range("H3").Value = "2016-03-11 09:19:44"
